$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1912.375
$ws.Range("J32").Value = 1912.375
$ws.Range("L32").Value = 1912.375
$ws.Range("N32").Value = -2564.375
$ws.Range("H69").Value = 4433.1665
$ws.Range("I69").Value = 3300
$ws.Range("J69").Value = 4999.75
$ws.Range("K69").Value = 9900
$ws.Range("L69").Value = 14999.25
$ws.Range("M69").Value = -9026
$ws.Range("N69").Value = -16747.25
$ws.Range("H72").Value = 4433.1665
$ws.Range("I72").Value = 3300
$ws.Range("J72").Value = 4999.75
$ws.Range("K72").Value = 29700
$ws.Range("L72").Value = 44997.75
$ws.Range("M72").Value = -25332
$ws.Range("N72").Value = -53733.75
$ws.Range("H87").Value = 25000
$ws.Range("J87").Value = 25000
$ws.Range("L87").Value = 25000
$ws.Range("N87").Value = -27496
$ws.Range("H90").Value = 25000
$ws.Range("J90").Value = 25000
$ws.Range("L90").Value = 75000
$ws.Range("N90").Value = -87480
$ws.Range("H129").Value = 912.27
$ws.Range("J129").Value = 929.2782999999999
$ws.Range("L129").Value = 2787.8349
$ws.Range("N129").Value = -12787.8349
$ws.Range("H137").Value = 2173.0715
$ws.Range("I137").Value = 1403.7222
$ws.Range("J137").Value = 2750.0833
$ws.Range("K137").Value = 4211.1666
$ws.Range("L137").Value = 8250.249899999999
$ws.Range("M137").Value = -1661.1666
$ws.Range("N137").Value = -13350.2499

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 4011.3076
$ws.Range("I3").Value = 4011.3076
$ws.Range("K3").Value = 4011.3076
$ws.Range("M3").Value = -3896.3076
$ws.Range("H32").Value = 8705.592000000001
$ws.Range("I32").Value = 8860.753000000001
$ws.Range("K32").Value = 8860.753000000001
$ws.Range("M32").Value = -8573.753000000001
$ws.Range("H132").Value = 5913.788
$ws.Range("I132").Value = 7252.15
$ws.Range("J132").Value = 3854.7693
$ws.Range("K132").Value = 21756.45
$ws.Range("L132").Value = 11564.3079
$ws.Range("M132").Value = -19226.45
$ws.Range("N132").Value = -16624.3079

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3499.5
$ws.Range("I105").Value = 3499.3333
$ws.Range("J105").Value = 3499.6667
$ws.Range("K105").Value = 3499.3333
$ws.Range("L105").Value = 3499.6667
$ws.Range("M105").Value = -1752.3333
$ws.Range("N105").Value = -6993.6667
$ws.Range("H134").Value = 2557.2424
$ws.Range("I134").Value = 2225
$ws.Range("J134").Value = 3443.2222
$ws.Range("K134").Value = 6675
$ws.Range("L134").Value = 10329.6666
$ws.Range("M134").Value = -4140
$ws.Range("N134").Value = -15399.6666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 52000
$ws.Range("J3").Value = 54000
$ws.Range("L3").Value = 54000
$ws.Range("N3").Value = -54226
$ws.Range("H99").Value = 1810
$ws.Range("J99").Value = 1691.6666
$ws.Range("L99").Value = 1691.6666
$ws.Range("N99").Value = -4687.6666
$ws.Range("H122").Value = 3164.3684
$ws.Range("I122").Value = 3774.5334
$ws.Range("J122").Value = 876.25
$ws.Range("K122").Value = 11323.6002
$ws.Range("L122").Value = 2628.75
$ws.Range("M122").Value = -8873.600199999999
$ws.Range("N122").Value = -7528.75
$ws.Range("H126").Value = 1810
$ws.Range("J126").Value = 1691.6666
$ws.Range("L126").Value = 5074.9998
$ws.Range("N126").Value = -10014.9998
$ws.Range("H132").Value = 1355693.5
$ws.Range("I132").Value = 2706087.5
$ws.Range("J132").Value = 5299.6
$ws.Range("K132").Value = 8118262.5
$ws.Range("L132").Value = 15898.8
$ws.Range("M132").Value = -8115732.5
$ws.Range("N132").Value = -20958.8
$ws.Range("H134").Value = 2147.3928
$ws.Range("I134").Value = 1194.95
$ws.Range("J134").Value = 4528.5
$ws.Range("K134").Value = 3584.85
$ws.Range("L134").Value = 13585.5
$ws.Range("M134").Value = -1049.85
$ws.Range("N134").Value = -18655.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 751.5
$ws.Range("J32").Value = 751.5
$ws.Range("L32").Value = 2254.5
$ws.Range("N32").Value = -2820.5
$ws.Range("H68").Value = 733.67706
$ws.Range("I68").Value = 603.9385
$ws.Range("J68").Value = 1005.70966
$ws.Range("K68").Value = 1811.8155
$ws.Range("L68").Value = 3017.12898
$ws.Range("M68").Value = -1000.8155
$ws.Range("N68").Value = -4639.12898
$ws.Range("H71").Value = 733.67706
$ws.Range("I71").Value = 603.9385
$ws.Range("J71").Value = 1005.70966
$ws.Range("K71").Value = 5435.4465
$ws.Range("L71").Value = 9051.38694
$ws.Range("M71").Value = -1379.4465
$ws.Range("N71").Value = -17163.38694
$ws.Range("H101").Value = 11944.214
$ws.Range("J101").Value = 12632.23
$ws.Range("L101").Value = 37896.69
$ws.Range("N101").Value = -42764.69
$ws.Range("H122").Value = 429.86667
$ws.Range("J122").Value = 319.64706
$ws.Range("L122").Value = 2876.82354
$ws.Range("N122").Value = -7776.82354
$ws.Range("H127").Value = 805.5
$ws.Range("J127").Value = 805.5
$ws.Range("L127").Value = 2416.5
$ws.Range("N127").Value = -12336.5
$ws.Range("H132").Value = 1074.5151
$ws.Range("I132").Value = 933.2222
$ws.Range("J132").Value = 1127.5
$ws.Range("K132").Value = 8398.9998
$ws.Range("L132").Value = 10147.5
$ws.Range("M132").Value = -5868.9998
$ws.Range("N132").Value = -15207.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3385.3333
$ws.Range("I132").Value = 2920.3635
$ws.Range("K132").Value = 8761.0905
$ws.Range("M132").Value = -6231.0905

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5116.1377
$ws.Range("I132").Value = 5285.4346
$ws.Range("J132").Value = 4467.1665
$ws.Range("K132").Value = 15856.3038
$ws.Range("L132").Value = 13401.4995
$ws.Range("M132").Value = -13326.3038
$ws.Range("N132").Value = -18461.4995
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 700.2
$ws.Range("I107").Value = 451
$ws.Range("J107").Value = 866.3333
$ws.Range("K107").Value = 1353
$ws.Range("L107").Value = 2598.9999
$ws.Range("M107").Value = 567
$ws.Range("N107").Value = -6438.9999
$ws.Range("H132").Value = 2558.1724
$ws.Range("I132").Value = 2446.2632
$ws.Range("J132").Value = 2770.8
$ws.Range("K132").Value = 7338.7896
$ws.Range("L132").Value = 8312.400000000001
$ws.Range("M132").Value = -4808.7896
$ws.Range("N132").Value = -13372.4
